$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: remove the old "_GoBack" bookmark (currently sitting at the
# end of the "User Researcher / ..." paragraph). It will be re-added
# inside the rewritten Solution paragraph below.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Step 2: delete everything from "Key Features" through the end of the
# document (Key Features, Process, Sprint 1-5, trailing empty para).
# ------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(13)
$pEnd = $d.Paragraphs.Item($d.Paragraphs.Count)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# Step 3: rewrite the "Solution" body paragraph (currently paragraph
# 10) with the new merged copy, and drop the "_GoBack" bookmark back
# in right after "... to buy and sell their".
# ------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$pStartPos = $p10.Range.Start
$pEndPos = $p10.Range.End
$solutionRange = $d.Range($pStartPos, $pEndPos)

$beforeBookmark = "Our team of experts created Encore Performance Gear, an online marketplace for the directors and designers of these ensembles to buy and sell their"
$afterBookmark = " used equipment. We used a user-centered agile development process that incorporated UX methodologies into five two-week development sprints. Our team deployed a minimal viable product to Heroku using React with Redux and a RESTful Django API."

$solutionRange.Text = $beforeBookmark + $afterBookmark

$bmPos = $pStartPos + $beforeBookmark.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# Step 4: insert the two new "UX Problem" heading + body paragraph
# pairs right after the Solution paragraph (now paragraph 10) and
# before "Value Proposition".
# ------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertParagraphAfter() | Out-Null
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = "UX Problem 1: Enable Searching & Filtering"
$p11.Style = "Heading 2"

$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertParagraphAfter() | Out-Null
$p12 = $d.Paragraphs.Item(12)
$p12.Range.Text = "Facebook Marketplace was difficult to search and navigate, Encore improved that by letting users search and filter through listings."
$p12.Style = "Normal"

$p12 = $d.Paragraphs.Item(12)
$p12.Range.InsertParagraphAfter() | Out-Null
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Text = "UX Problem 2: Building Trust on the Platform"
$p13.Style = "Heading 2"

$p13 = $d.Paragraphs.Item(13)
$p13.Range.InsertParagraphAfter() | Out-Null
$p14 = $d.Paragraphs.Item(14)
$p14.Range.Text = "Our initial user research taught us that trust was a large factor with these large purchases. We wanted to include profiles to help users feel better about who they" + [char]0x2019 + "re buying and selling to. "
$p14.Style = "Normal"

Write-Host "Done. Final paragraph count: $($d.Paragraphs.Count)"
